# Append the new task rows (rows 6-27) to the "Tasks" sheet and grow the
# used range from A1:F5 to A1:F27, matching the uploaded-workbook diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# $null below means "leave the cell blank, but still typed as an explicit
# empty string" (mirrors the source file's <c t="str"><v/></c> cells).
$newRows = @(
    @{ Row=6;  A="pc room";  B="scvsvsvv";       C="medium"; D="samith";  E="In Progress" }
    @{ Row=7;  A="pc room";  B="plingcfxrtvdk";   C="low";    D="sumith";  E="In Progress" }
    @{ Row=8;  A="pc room";  B="plplpl";          C="low";    D="plo";     E="In Progress" }
    @{ Row=9;  A="pc room";  B="plpl`n";          C="low";    D="pathum";  E="In Progress" }
    @{ Row=10; A="plc room"; B="hbsgvsyvu";       C="low";    D="dulan";   E="In Progress" }
    @{ Row=11; A="pc room";  B="pooipp";          C="medium"; D="pl";      E="In Progress" }
    @{ Row=12; A="pc room";  B="fgfb";            C="low";    D="plpl";    E="In Progress" }
    @{ Row=13; A="grgb";     B="fbb";             C="low";    D="dd";      E="In Progress" }
    @{ Row=14; A="grgrebfb"; B="bfbbf";           C="low";    D=$null;     E="Pending" }
    @{ Row=15; A="abfbf";    B="ab";              C="low";    D=$null;     E="Pending" }
    @{ Row=16; A="agrgrrb";  B="aagg";            C="low";    D="pathum";  E="In Progress" }
    @{ Row=17; A="agrg";     B="agr";             C="low";    D=$null;     E="Pending" }
    @{ Row=18; A="egg";      B="sgtg";            C="low";    D=$null;     E="Pending" }
    @{ Row=19; A="rgrb";     B="bbb";             C="low";    D="yy";      E="In Progress" }
    @{ Row=20; A="plc room"; B="ll";              C="low";    D="hhh";     E="In Progress" }
    @{ Row=21; A="pc room";  B="cc";              C="low";    D="vv";      E="In Progress" }
    @{ Row=22; A="pc room";  B="ppp";             C="medium"; D="sumith";  E="In Progress" }
    @{ Row=23; A="plc";      B="fggk";            C="medium"; D="samith";  E="In Progress" }
    @{ Row=24; A="dd";       B="ddd";             C="low";    D=$null;     E="Pending" }
    @{ Row=25; A="pc room";  B="cc";              C="medium"; D=$null;     E="Pending" }
    @{ Row=26; A="pc room";  B="kkk";             C="low";    D="hgj";     E="In Progress" }
    @{ Row=27; A="dd";       B="d";               C="low";    D=$null;     E="Pending" }
)

function Set-EmptyTextCell($cell) {
    # A plain `.Value = ""` leaves the cell completely unset (no <c> element
    # at all once saved), but the target file has an explicit empty *text*
    # cell there. Entering a lone leading apostrophe is the standard Excel
    # way to force text-type on an empty cell; resetting the style back to
    # "Normal" afterwards drops the quote-prefix formatting it introduces so
    # the cell ends up as plain empty text with the default style.
    $cell.Value = "'"
    $cell.Style = "Normal"
}

foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C

    if ($null -eq $r.D) {
        Set-EmptyTextCell $ws.Range("D$($r.Row)")
    } else {
        $ws.Range("D$($r.Row)").Value = $r.D
    }

    $ws.Range("E$($r.Row)").Value = $r.E

    # "Time Taken (mins)" column is always blank (explicit empty text) for
    # every row, old and new alike.
    Set-EmptyTextCell $ws.Range("F$($r.Row)")
}

Write-Output "Added $($newRows.Count) rows to $($ws.Name)"
